# Refresh the crypto price/volume snapshot in the worksheet.
# Mirrors the upstream GitHub Actions scraper commit:
#   "Updated cryptos list on Tue Oct  3 02:34:01 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.552.25"
$ws.Range("D3").Value = "1.666.73"
$ws.Range("E3").Value = "  -3.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.21"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.51"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "1.901.95"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "1.677.94"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "249.91"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "27.570.71"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("E24").Value = "  -5.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.58"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.55"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +3.83%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").Value = "1.475.76"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  -5.48%  "
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.942"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -5.84%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.76"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E43").Value = "  -6.93%  "
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "1.810.16"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.46"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "41.97"
$ws.Range("E50").Value = "  +16.01%  "
$ws.Range("E51").Value = "  -3.32%  "
